$wb = $excel.ActiveWorkbook
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws8 = $wb.Worksheets.Add($null, $last)
$ws8.Name = "Credit rating correlation"

$ws8.Range("A1").Value = "Rating"
$ws8.Range("B1").Value = "Revenue"
$ws8.Range("C1").Value = "Debt Ratio"

$ws8.Range("A2").Value = 0
$ws8.Range("B2").Value = 5799.0640000000003
$ws8.Range("C2").Value = 418.714

$ws8.Range("A3").Value = 1
$ws8.Range("B3").Value = 6893.4979999999996
$ws8.Range("C3").Value = 174.69800000000001

$ws8.Range("A4").Value = 2
$ws8.Range("B4").Value = 7198.6459999999997
$ws8.Range("C4").Value = 263.03100000000001

$ws8.Range("A5").Value = 3
$ws8.Range("B5").Value = 7490.268
$ws8.Range("C5").Value = 227.625

$ws8.Range("A6").Value = 4
$ws8.Range("B6").Value = 8302.8709999999992
$ws8.Range("C6").Value = 181.32400000000001

$ws8.Range("A7").Value = 5
$ws8.Range("B7").Value = 8587.1630000000005
$ws8.Range("C7").Value = 172.709

$ws8.Range("A8").Value = 6
$ws8.Range("B8").Value = 8446.1540000000005
$ws8.Range("C8").Value = 179.23

$ws8.Range("A9").Value = 7
$ws8.Range("B9").Value = 8536.2309999999998
$ws8.Range("C9").Value = 300.82799999999997

$ws8.Range("A10").Value = 8
$ws8.Range("B10").Value = 4910
$ws8.Range("C10").Value = 0.20599999999999999

$ws8.Range("A11").Value = 9
$ws8.Range("B11").Value = 2789
$ws8.Range("C11").Value = 0.64400000000000002

$ws8.Range("A12").Value = 10
$ws8.Range("B12").Value = 7500
$ws8.Range("C12").Value = 0.13400000000000001

$ws8.Range("E2").Value = "Revenue corr."
$ws8.Range("F2").Formula = "=CORREL(A2:A8,B2:B8)"

$ws8.Range("E3").Value = "Debt ratio corr."
$ws8.Range("F3").Formula = "=CORREL(A2:A8,C2:C8)"

$ws8.Columns.AutoFit()
